# Rename sheet "adj_r_squared" -> "metrics" and populate it with the
# Adj.R^2 / NRMSE / SMAPE / RMSE metrics table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("adj_r_squared")
$ws.Name = "metrics"

$ws.Range("A1").Value = "metrics"
$ws.Range("B1").Value = "value"

$ws.Range("A2").Value = "Adj.R^2"
$ws.Range("B2").Value = 0.665451304382641

$ws.Range("A3").Value = "NRMSE"
$ws.Range("B3").Value = 0.02097291346767

$ws.Range("A4").Value = "SMAPE"
$ws.Range("B4").Value = 0.0163201283064115

$ws.Range("A5").Value = "RMSE"
$ws.Range("B5").Value = 0.213571634585396
